# Update cryptocurrency price/volume data to match the latest GitHub Actions scrape.
# (commit: "Updated cryptos list on Mon Sep  9 05:12:22 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "54.560.16"
$ws.Range("E2").Value = "  +0.28%  "
# Row 3
$ws.Range("D3").Value = "2.276.72"
$ws.Range("E3").Value = "  -0.24%  "
# Row 4
$ws.Range("E4").Value = "  +0.08%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D5").Value = "503.22"
$ws.Range("E5").Value = "  +1.00%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D6").Value = "128.44"
$ws.Range("E6").Value = "  +0.06%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.13%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  -0.32%  "
# Row 9
$ws.Range("D9").Value = "2.293.01"
$ws.Range("E9").Value = "  -0.15%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D10").Value = "0.0964"
$ws.Range("E10").Value = "  +0.57%  "
# Row 11
$ws.Range("E11").Value = "  +0.70%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  +2.67%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D13").Value = "4.88"
$ws.Range("E13").Value = "  +3.69%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D14").Value = "23.29"
$ws.Range("E14").Value = "  +3.68%  "
# Row 15
$ws.Range("D15").Value = "2.682.49"
$ws.Range("E15").Value = "  -0.25%  "
# Row 16
$ws.Range("D16").Value = "54.640.45"
$ws.Range("E16").Value = "  +0.56%  "
# Row 17
$ws.Range("E17").Value = "  +0.67%  "
# Row 18
$ws.Range("D18").Value = "2.286.81"
$ws.Range("E18").Value = "  -1.18%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D19").Value = "10.31"
$ws.Range("E19").Value = "  +0.64%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D20").Value = "4.12"
$ws.Range("E20").Value = "  -0.46%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D21").Value = "306.45"
$ws.Range("E21").Value = "  +0.17%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D22").Value = "6.42"
$ws.Range("E22").Value = "  -0.44%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.11%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D24").Value = "60.15"
$ws.Range("E24").Value = "  -3.17%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D25").Value = "0.995"
$ws.Range("E25").Value = "  -0.35%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D26").Value = "0.150"
$ws.Range("E26").Value = "  -0.61%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D27").Value = "7.41"
$ws.Range("E27").Value = "  +1.09%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D28").Value = "170.77"
$ws.Range("E28").Value = "  -1.75%  "
# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D29").Value = "1.62"
$ws.Range("E29").Value = "  +0.08%  "
# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0698"
$ws.Range("E30").Value = "  +1.22%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D31").Value = "6.03"
$ws.Range("E31").Value = "  +1.21%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  +1.95%  "
# Row 33
$ws.Range("E33").Value = "  +0.03%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D34").Value = "17.91"
$ws.Range("E34").Value = "  +0.56%  "
# Row 35
$ws.Range("E35").Value = "  -0.11%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D36").Value = "0.909"
$ws.Range("E36").Value = "  -1.81%  "
# Row 37
$ws.Range("E37").Value = "  -0.77%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D38").Value = "3.77"
$ws.Range("E38").Value = "  +0.07%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D39").Value = "36.39"
$ws.Range("E39").Value = "  +1.41%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D40").Value = "0.376"
$ws.Range("E40").Value = "  +0.42%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D41").Value = "1.41"
$ws.Range("E41").Value = "  -0.32%  "
# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D42").Value = "5.04"
$ws.Range("E42").Value = "  +5.03%  "
# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  -0.33%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D44").Value = "126.45"
$ws.Range("E44").Value = "  +0.12%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D45").Value = "249.32"
$ws.Range("E45").Value = "  +3.65%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D46").Value = "0.0496"
$ws.Range("E46").Value = "  +0.71%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D47").Value = "0.0899"
$ws.Range("E47").Value = "  +0.08%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D48").Value = "0.547"
$ws.Range("E48").Value = "  -0.27%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"   # keep as text, like the source data
$ws.Range("D49").Value = "0.374"
$ws.Range("E49").Value = "  +0.46%  "
# Row 50
$ws.Range("E50").Value = "  -0.21%  "
# Row 51
$ws.Range("E51").Value = "  +0.41%  "

# Restore default (unformatted) cell style now that the text values are committed,
# so number formatting matches the original workbook.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
